$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove any pre-existing formatting on the whole used range so the
# leftover "asignados" number style (s="1") isn't dragged onto the
# reshuffled columns.
$ws.Cells.ClearFormats() | Out-Null

# ---- Header row -----------------------------------------------------
$ws.Range("A1").Value = "codigo"
$ws.Range("B1").Value = "nombre"
$ws.Range("C1").Value = "semestre"
$ws.Range("D1").Value = "asignados"
$ws.Range("E1").Value = "no_periodos"
$ws.Range("F1").Value = "carrera"
$ws.Range("G1").Value = "area"

# ---- Data rows (codigo, nombre, semestre, asignados, no_periodos, carrera, area)
$row = 2
$ws.Range("A$row").Value = "M001"
$ws.Range("B$row").Value = "Materia 01"
$ws.Range("C$row").Value = "I"
$ws.Range("D$row").Value = 50
$ws.Range("E$row").Value = 1
$ws.Range("F$row").Value = "C001"
$ws.Range("G$row").Value = "BASICA"

$row = 3
$ws.Range("A$row").Value = "M002"
$ws.Range("B$row").Value = "Materia 02"
$ws.Range("C$row").Value = "V"
$ws.Range("D$row").Value = 25
$ws.Range("E$row").Value = 1
$ws.Range("F$row").Value = "C002"
$ws.Range("G$row").Value = "PROFESIONAL"

$row = 4
$ws.Range("A$row").Value = "M003"
$ws.Range("B$row").Value = "Materia 03"
$ws.Range("C$row").Value = "II"
$ws.Range("D$row").Value = 30
$ws.Range("E$row").Value = 1
$ws.Range("F$row").Value = "C003"
$ws.Range("G$row").Value = "BASICA"

$row = 5
$ws.Range("A$row").Value = "M004"
$ws.Range("B$row").Value = "Materia 04"
$ws.Range("C$row").Value = "I"
$ws.Range("D$row").Value = 40
$ws.Range("E$row").Value = 1
$ws.Range("F$row").Value = "C004"
$ws.Range("G$row").Value = "BASICA"

$row = 6
$ws.Range("A$row").Value = "M005"
$ws.Range("B$row").Value = "Materia 05"
$ws.Range("C$row").Value = "III"
$ws.Range("D$row").Value = 40
$ws.Range("E$row").Value = 1
$ws.Range("F$row").Value = "C005"
$ws.Range("G$row").Value = "BASICA"

$row = 7
$ws.Range("A$row").Value = "M006"
$ws.Range("B$row").Value = "Materia 06"
$ws.Range("C$row").Value = "II"
$ws.Range("D$row").Value = 75
$ws.Range("E$row").Value = 1
$ws.Range("F$row").Value = "C006"
$ws.Range("G$row").Value = "BASICA"

$row = 8
$ws.Range("A$row").Value = "M007"
$ws.Range("B$row").Value = "Materia 07"
$ws.Range("C$row").Value = "IV"
$ws.Range("D$row").Value = 100
$ws.Range("E$row").Value = 1
$ws.Range("F$row").Value = "C001"
$ws.Range("G$row").Value = "PROFESIONAL"

$row = 9
$ws.Range("A$row").Value = "M008"
$ws.Range("B$row").Value = "Materia 08"
$ws.Range("C$row").Value = "I"
$ws.Range("D$row").Value = 125
$ws.Range("E$row").Value = 1
$ws.Range("F$row").Value = "C002"
$ws.Range("G$row").Value = "BASICA"

$row = 10
$ws.Range("A$row").Value = "M009"
$ws.Range("B$row").Value = "Materia 09"
$ws.Range("C$row").Value = "III"
$ws.Range("D$row").Value = 80
$ws.Range("E$row").Value = 1
$ws.Range("F$row").Value = "C003"
$ws.Range("G$row").Value = "PROFESIONAL"

$row = 11
$ws.Range("A$row").Value = "M010"
$ws.Range("B$row").Value = "Materia 10"
$ws.Range("C$row").Value = "V"
$ws.Range("D$row").Value = 30
$ws.Range("E$row").Value = 1
$ws.Range("F$row").Value = "C004"
$ws.Range("G$row").Value = "PROFESIONAL"

$row = 12
$ws.Range("A$row").Value = "M011"
$ws.Range("B$row").Value = "Materia 11"
$ws.Range("C$row").Value = "IV"
$ws.Range("D$row").Value = 20
$ws.Range("E$row").Value = 1
$ws.Range("F$row").Value = "C005"
$ws.Range("G$row").Value = "PROFESIONAL"

$row = 13
$ws.Range("A$row").Value = "M012"
$ws.Range("B$row").Value = "Inter 1"
$ws.Range("C$row").Value = "II"
$ws.Range("D$row").Value = 15
$ws.Range("E$row").Value = 2
$ws.Range("F$row").Value = "C006"
$ws.Range("G$row").Value = "BASICA"

# ---- Column widths (best-fit driven by the new "no_periodos" / "area"
# content, as Excel computed when the author widened those columns) ----
$ws.Columns("E:E").ColumnWidth = 11.33
$ws.Columns("G:G").ColumnWidth = 12.45

# ---- Sheet view / selection -------------------------------------------
$ws.Range("G17").Select() | Out-Null

# ---- Page setup ---------------------------------------------------------
$ws.PageSetup.Orientation = 1

Write-Output "ok"
